$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string (e.g. "398.74") must be
# forced to Text format first, otherwise Excel auto-converts them to real
# numbers on assignment and the original text formatting (trailing zeros,
# thousand-dot grouping, etc.) would be lost.
$textForcedUpdates = @(
  @{Row=5; Col=4; Value="398.74"},
  @{Row=6; Col=4; Value="109.00"},
  @{Row=7; Col=4; Value="0.579"},
  @{Row=9; Col=4; Value="0.621"},
  @{Row=10; Col=4; Value="39.31"},
  @{Row=14; Col=4; Value="8.29"},
  @{Row=15; Col=4; Value="18.98"},
  @{Row=18; Col=4; Value="11.10"},
  @{Row=20; Col=4; Value="3.32"},
  @{Row=21; Col=4; Value="0.0000107"},
  @{Row=23; Col=4; Value="296.44"},
  @{Row=24; Col=4; Value="74.32"},
  @{Row=26; Col=4; Value="28.09"},
  @{Row=27; Col=4; Value="4.40"},
  @{Row=28; Col=4; Value="7.89"},
  @{Row=29; Col=4; Value="7.41"},
  @{Row=33; Col=4; Value="11.25"},
  @{Row=34; Col=4; Value="40.43"},
  @{Row=35; Col=4; Value="0.0500"},
  @{Row=37; Col=4; Value="51.83"},
  @{Row=38; Col=4; Value="3.11"},
  @{Row=39; Col=4; Value="0.999"},
  @{Row=40; Col=4; Value="3.48"},
  @{Row=41; Col=4; Value="138.09"},
  @{Row=43; Col=4; Value="0.284"},
  @{Row=45; Col=4; Value="16.84"},
  @{Row=46; Col=4; Value="3.90"},
  @{Row=47; Col=4; Value="22.32"},
  @{Row=48; Col=4; Value="2.21"},
  @{Row=50; Col=4; Value="2.47"},
  @{Row=51; Col=4; Value="1.92"}
)

foreach ($item in $textForcedUpdates) {
  $cell = $ws.Cells.Item($item.Row, $item.Col)
  $cell.NumberFormat = "@"
  $cell.Value2 = $item.Value
}

# Remaining changed cells (plain text: names, links, and the padded
# percentage strings in column E) can be written directly.
$plainUpdates = @(
  @{Row=2; Col=4; Value="57.544.46"},
  @{Row=2; Col=5; Value="  +2.95%  "},
  @{Row=3; Col=4; Value="3.270.62"},
  @{Row=3; Col=5; Value="  +1.69%  "},
  @{Row=4; Col=5; Value="  -0.10%  "},
  @{Row=5; Col=5; Value="  +0.82%  "},
  @{Row=6; Col=5; Value="  -1.30%  "},
  @{Row=7; Col=5; Value="  +5.18%  "},
  @{Row=8; Col=5; Value="  -0.06%  "},
  @{Row=9; Col=5; Value="  +0.77%  "},
  @{Row=10; Col=5; Value="  +0.56%  "},
  @{Row=11; Col=5; Value="  +5.97%  "},
  @{Row=12; Col=5; Value="  +1.27%  "},
  @{Row=13; Col=4; Value="3.786.77"},
  @{Row=13; Col=5; Value="  +1.48%  "},
  @{Row=14; Col=5; Value="  +3.11%  "},
  @{Row=15; Col=5; Value="  +0.02%  "},
  @{Row=16; Col=4; Value="3.256.37"},
  @{Row=16; Col=5; Value="  +1.25%  "},
  @{Row=17; Col=5; Value="  -0.95%  "},
  @{Row=18; Col=5; Value="  +2.54%  "},
  @{Row=19; Col=4; Value="57.291.43"},
  @{Row=19; Col=5; Value="  +2.62%  "},
  @{Row=20; Col=5; Value="  +0.18%  "},
  @{Row=21; Col=5; Value="  +5.39%  "},
  @{Row=22; Col=5; Value="  -0.01%  "},
  @{Row=23; Col=5; Value="  -0.51%  "},
  @{Row=25; Col=5; Value="  -0.53%  "},
  @{Row=26; Col=5; Value="  +0.20%  "},
  @{Row=27; Col=2; Value="LEO"},
  @{Row=27; Col=3; Value="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"},
  @{Row=27; Col=5; Value="  +0.83%  "},
  @{Row=28; Col=2; Value="Filecoin"},
  @{Row=28; Col=3; Value="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"},
  @{Row=28; Col=5; Value="  -2.93%  "},
  @{Row=29; Col=5; Value="  -1.01%  "},
  @{Row=30; Col=5; Value="  -1.29%  "},
  @{Row=31; Col=5; Value="  -0.42%  "},
  @{Row=32; Col=5; Value="  +1.64%  "},
  @{Row=33; Col=5; Value="  +1.53%  "},
  @{Row=34; Col=5; Value="  +12.46%  "},
  @{Row=35; Col=5; Value="  +2.29%  "},
  @{Row=36; Col=5; Value="  +0.64%  "},
  @{Row=37; Col=5; Value="  +0.94%  "},
  @{Row=38; Col=2; Value="Stacks"},
  @{Row=38; Col=3; Value="https://coinranking.com/coin/mMPrMcB7+stacks-stx"},
  @{Row=38; Col=5; Value="  +0.20%  "},
  @{Row=39; Col=2; Value="FirstDigitalUSD"},
  @{Row=39; Col=3; Value="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"},
  @{Row=39; Col=5; Value="  -0.09%  "},
  @{Row=40; Col=5; Value="  -1.13%  "},
  @{Row=41; Col=5; Value="  +1.78%  "},
  @{Row=42; Col=5; Value="  +2.13%  "},
  @{Row=43; Col=5; Value="  +0.39%  "},
  @{Row=44; Col=5; Value="  -1.75%  "},
  @{Row=45; Col=5; Value="  -2.55%  "},
  @{Row=46; Col=5; Value="  -1.83%  "},
  @{Row=47; Col=5; Value="  +0.98%  "},
  @{Row=48; Col=5; Value="  +4.42%  "},
  @{Row=49; Col=4; Value="2.157.94"},
  @{Row=49; Col=5; Value="  +1.48%  "},
  @{Row=50; Col=5; Value="  -0.05%  "},
  @{Row=51; Col=5; Value="  -11.26%  "}
)

foreach ($item in $plainUpdates) {
  $ws.Cells.Item($item.Row, $item.Col).Value2 = $item.Value
}

Write-Output "Applied $($textForcedUpdates.Count + $plainUpdates.Count) cell updates"
